$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: add end time (copy time format from B4) and accomplishment text
$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = 0.31041666666666667
$ws.Range("D4").Value = "Finished putting in procedural midi plugin, looked into in-engine synthesis, recorded a basic scale and drum beats"

# Row 5: new work log entry
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 43729

$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = 0.33958333333333335

$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = 0.4597222222222222

$ws.Range("D5").Value = "Changed system for spawning grid tiles from editor, made it so that default towers can be spawned from grid tiles in editor"

$excel.CutCopyMode = $false

$ws.Range("D5").Select()
